$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price-report row was added to the data set. In the sheet it
# lands at row 63, pushing the former rows 63:81 down to 64:82 (the last
# existing row, old row 81, ends up duplicated as the new final row 82).
$ws.Rows("63:63").Insert()

# Populate the newly inserted row 63 with the new observation.
$ws.Range("A63").Value = 8
$ws.Range("B63").Value = "Terminal La Palmera de La Serena"
$ws.Range("C63").Value = "Coquimbo"
$ws.Range("D63").Value2 = 44463
$ws.Range("E63").Value = 4
$ws.Range("F63").Value = 100112001
$ws.Range("G63").Value = "Berenjena"
$ws.Range("H63").Value = "Sin especificar"
$ws.Range("I63").Value = "Primera"
$ws.Range("J63").Value = 600
$ws.Range("K63").Value = 9000
$ws.Range("L63").Value = 10000
$ws.Range("M63").Value = 9500
$ws.Range("N63").Value = "$/caja 60 unidades"
$ws.Range("O63").Value = "Región de Arica y Parinacota"
$ws.Range("P63").Value = 158
$ws.Range("Q63").Value = 60
$ws.Range("R63").Value = "Hortaliza"
